# Documentation Update - April 28th
# Added documentation for:
# - Header of ActionKeywordsFile
# - Documented Click and Navigate actions.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook-level view metadata (window size / OLE display size).
# The OLE/embed display area becomes A1:E13 on the first sheet.
# ---------------------------------------------------------------------
$wsSuite = $wb.Worksheets.Item("Test Suite")
$wsSuite.Activate()
$wsSuite.Range("A1:E13").Select()
$excel.ActiveWindow.WindowWidth = 16275

# ---------------------------------------------------------------------
# Sheet "Test Suite": flip C6 from "yes" to "no", move selection to C5
# ---------------------------------------------------------------------
$wsSuite.Range("C6").Value = "no"
$wsSuite.Range("C5").Select()

# ---------------------------------------------------------------------
# Sheet "Test Cases": the two test-case blocks (rows 2-11 and rows
# 14-26) swap their Test Case ID (column A) values with each other,
# scroll back to the top-left of the sheet, and move the selection.
# ---------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Activate()

for ($r = 2; $r -le 11; $r++) {
    $wsCases.Range("A$r").Value = "Validate Mortgage Page"
}

for ($r = 14; $r -le 26; $r++) {
    $wsCases.Range("A$r").Value = "Validate With Mahesh"
}

$wsCases.Range("A1").Select()
$wsCases.Range("B7").Select()

# ---------------------------------------------------------------------
# Sheet "Result": move selection to D5
# ---------------------------------------------------------------------
$wsResult = $wb.Worksheets.Item("Result")
$wsResult.Activate()
$wsResult.Range("D5").Select()

# ---------------------------------------------------------------------
# Workbook-level view metadata (window size / OLE display size).
# ---------------------------------------------------------------------
$wsSuite.Activate()
$wsSuite.Range("A1:E13").Select()
$excel.ActiveWindow.WindowWidth = 16275

# "Test Cases" stays the active/selected tab when the workbook is saved.
$wsCases.Activate()
$wsCases.Range("B7").Select()
